# "Generate Report for Handoff"
#
# The localization status report is regenerated. For the de-de locale, two
# files (rows 10 and 12, corresponding to the C10/C12 handoff filenames)
# were just re-handed-off, so their "Latest Handoff Datetime" (column D)
# advances to the new handoff timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("D10").Value = "2016-03-09 12:23:28"
$ws.Range("D12").Value = "2016-03-09 12:23:28"
